$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each cell is written via a leading-apostrophe ("quote-prefix") literal so
# Excel stores the exact text (preserving things like "41.628.41",
# "6.90"'s trailing zero, or the "  +0.17%  " padding) instead of silently
# reinterpreting numeric-looking strings as numbers/percentages. The cell's
# original Style is captured and reapplied afterwards so the quote-prefix
# flag Excel sets on the cell format does not leave a stray style/format
# change behind.

# Row 2
$origStyle = $ws.Cells.Item(2, 4).Style
$ws.Cells.Item(2, 4).Value = '''41.628.41'
$ws.Cells.Item(2, 4).Style = $origStyle
$origStyle = $ws.Cells.Item(2, 5).Style
$ws.Cells.Item(2, 5).Value = '''  +0.17%  '
$ws.Cells.Item(2, 5).Style = $origStyle

# Row 3
$origStyle = $ws.Cells.Item(3, 4).Style
$ws.Cells.Item(3, 4).Value = '''2.473.91'
$ws.Cells.Item(3, 4).Style = $origStyle
$origStyle = $ws.Cells.Item(3, 5).Style
$ws.Cells.Item(3, 5).Value = '''  +0.63%  '
$ws.Cells.Item(3, 5).Style = $origStyle

# Row 5
$origStyle = $ws.Cells.Item(5, 4).Style
$ws.Cells.Item(5, 4).Value = '''318.96'
$ws.Cells.Item(5, 4).Style = $origStyle
$origStyle = $ws.Cells.Item(5, 5).Style
$ws.Cells.Item(5, 5).Value = '''  +1.44%  '
$ws.Cells.Item(5, 5).Style = $origStyle

# Row 6
$origStyle = $ws.Cells.Item(6, 4).Style
$ws.Cells.Item(6, 4).Value = '''92.37'
$ws.Cells.Item(6, 4).Style = $origStyle
$origStyle = $ws.Cells.Item(6, 5).Style
$ws.Cells.Item(6, 5).Value = '''  +1.18%  '
$ws.Cells.Item(6, 5).Style = $origStyle

# Row 7
$origStyle = $ws.Cells.Item(7, 5).Style
$ws.Cells.Item(7, 5).Value = '''  +0.75%  '
$ws.Cells.Item(7, 5).Style = $origStyle

# Row 8
$origStyle = $ws.Cells.Item(8, 5).Style
$ws.Cells.Item(8, 5).Value = '''  +0.01%  '
$ws.Cells.Item(8, 5).Style = $origStyle

# Row 9
$origStyle = $ws.Cells.Item(9, 4).Style
$ws.Cells.Item(9, 4).Value = '''0.514'
$ws.Cells.Item(9, 4).Style = $origStyle
$origStyle = $ws.Cells.Item(9, 5).Style
$ws.Cells.Item(9, 5).Value = '''  +0.95%  '
$ws.Cells.Item(9, 5).Style = $origStyle

# Row 10
$origStyle = $ws.Cells.Item(10, 4).Style
$ws.Cells.Item(10, 4).Value = '''0.0867'
$ws.Cells.Item(10, 4).Style = $origStyle
$origStyle = $ws.Cells.Item(10, 5).Style
$ws.Cells.Item(10, 5).Value = '''  +8.89%  '
$ws.Cells.Item(10, 5).Style = $origStyle

# Row 11
$origStyle = $ws.Cells.Item(11, 4).Style
$ws.Cells.Item(11, 4).Value = '''33.15'
$ws.Cells.Item(11, 4).Style = $origStyle
$origStyle = $ws.Cells.Item(11, 5).Style
$ws.Cells.Item(11, 5).Value = '''  +2.03%  '
$ws.Cells.Item(11, 5).Style = $origStyle

# Row 12
$origStyle = $ws.Cells.Item(12, 5).Style
$ws.Cells.Item(12, 5).Value = '''  -0.04%  '
$ws.Cells.Item(12, 5).Style = $origStyle

# Row 13
$origStyle = $ws.Cells.Item(13, 4).Style
$ws.Cells.Item(13, 4).Value = '''2.854.73'
$ws.Cells.Item(13, 4).Style = $origStyle
$origStyle = $ws.Cells.Item(13, 5).Style
$ws.Cells.Item(13, 5).Value = '''  +0.58%  '
$ws.Cells.Item(13, 5).Style = $origStyle

# Row 14
$origStyle = $ws.Cells.Item(14, 4).Style
$ws.Cells.Item(14, 4).Value = '''6.90'
$ws.Cells.Item(14, 4).Style = $origStyle
$origStyle = $ws.Cells.Item(14, 5).Style
$ws.Cells.Item(14, 5).Value = '''  +1.01%  '
$ws.Cells.Item(14, 5).Style = $origStyle

# Row 15
$origStyle = $ws.Cells.Item(15, 4).Style
$ws.Cells.Item(15, 4).Value = '''15.49'
$ws.Cells.Item(15, 4).Style = $origStyle
$origStyle = $ws.Cells.Item(15, 5).Style
$ws.Cells.Item(15, 5).Value = '''  -1.64%  '
$ws.Cells.Item(15, 5).Style = $origStyle

# Row 16
$origStyle = $ws.Cells.Item(16, 4).Style
$ws.Cells.Item(16, 4).Value = '''2.463.42'
$ws.Cells.Item(16, 4).Style = $origStyle
$origStyle = $ws.Cells.Item(16, 5).Style
$ws.Cells.Item(16, 5).Value = '''  +0.67%  '
$ws.Cells.Item(16, 5).Style = $origStyle

# Row 17
$origStyle = $ws.Cells.Item(17, 5).Style
$ws.Cells.Item(17, 5).Value = '''  +2.85%  '
$ws.Cells.Item(17, 5).Style = $origStyle

# Row 18
$origStyle = $ws.Cells.Item(18, 4).Style
$ws.Cells.Item(18, 4).Value = '''41.589.43'
$ws.Cells.Item(18, 4).Style = $origStyle
$origStyle = $ws.Cells.Item(18, 5).Style
$ws.Cells.Item(18, 5).Value = '''  +0.08%  '
$ws.Cells.Item(18, 5).Style = $origStyle

# Row 19
$origStyle = $ws.Cells.Item(19, 4).Style
$ws.Cells.Item(19, 4).Value = '''6.45'
$ws.Cells.Item(19, 4).Style = $origStyle
$origStyle = $ws.Cells.Item(19, 5).Style
$ws.Cells.Item(19, 5).Value = '''  -0.29%  '
$ws.Cells.Item(19, 5).Style = $origStyle

# Row 20
$origStyle = $ws.Cells.Item(20, 4).Style
$ws.Cells.Item(20, 4).Value = '''0.0₃0945'
$ws.Cells.Item(20, 4).Style = $origStyle
$origStyle = $ws.Cells.Item(20, 5).Style
$ws.Cells.Item(20, 5).Value = '''  +0.94%  '
$ws.Cells.Item(20, 5).Style = $origStyle

# Row 21
$origStyle = $ws.Cells.Item(21, 4).Style
$ws.Cells.Item(21, 4).Value = '''70.69'
$ws.Cells.Item(21, 4).Style = $origStyle
$origStyle = $ws.Cells.Item(21, 5).Style
$ws.Cells.Item(21, 5).Value = '''  -0.25%  '
$ws.Cells.Item(21, 5).Style = $origStyle

# Row 22
$origStyle = $ws.Cells.Item(22, 4).Style
$ws.Cells.Item(22, 4).Value = '''11.28'
$ws.Cells.Item(22, 4).Style = $origStyle
$origStyle = $ws.Cells.Item(22, 5).Style
$ws.Cells.Item(22, 5).Value = '''  -0.31%  '
$ws.Cells.Item(22, 5).Style = $origStyle

# Row 23
$origStyle = $ws.Cells.Item(23, 4).Style
$ws.Cells.Item(23, 4).Value = '''240.47'
$ws.Cells.Item(23, 4).Style = $origStyle
$origStyle = $ws.Cells.Item(23, 5).Style
$ws.Cells.Item(23, 5).Value = '''  +1.53%  '
$ws.Cells.Item(23, 5).Style = $origStyle

# Row 24
$origStyle = $ws.Cells.Item(24, 5).Style
$ws.Cells.Item(24, 5).Value = '''  +1.57%  '
$ws.Cells.Item(24, 5).Style = $origStyle

# Row 25
$origStyle = $ws.Cells.Item(25, 4).Style
$ws.Cells.Item(25, 4).Value = '''1.95'
$ws.Cells.Item(25, 4).Style = $origStyle
$origStyle = $ws.Cells.Item(25, 5).Style
$ws.Cells.Item(25, 5).Value = '''  +2.84%  '
$ws.Cells.Item(25, 5).Style = $origStyle

# Row 27
$origStyle = $ws.Cells.Item(27, 4).Style
$ws.Cells.Item(27, 4).Value = '''24.80'
$ws.Cells.Item(27, 4).Style = $origStyle
$origStyle = $ws.Cells.Item(27, 5).Style
$ws.Cells.Item(27, 5).Value = '''  +2.23%  '
$ws.Cells.Item(27, 5).Style = $origStyle

# Row 28
$origStyle = $ws.Cells.Item(28, 4).Style
$ws.Cells.Item(28, 4).Value = '''2.23'
$ws.Cells.Item(28, 4).Style = $origStyle
$origStyle = $ws.Cells.Item(28, 5).Style
$ws.Cells.Item(28, 5).Value = '''  -1.51%  '
$ws.Cells.Item(28, 5).Style = $origStyle

# Row 29
$origStyle = $ws.Cells.Item(29, 4).Style
$ws.Cells.Item(29, 4).Value = '''9.70'
$ws.Cells.Item(29, 4).Style = $origStyle
$origStyle = $ws.Cells.Item(29, 5).Style
$ws.Cells.Item(29, 5).Value = '''  +0.69%  '
$ws.Cells.Item(29, 5).Style = $origStyle

# Row 30
$origStyle = $ws.Cells.Item(30, 4).Style
$ws.Cells.Item(30, 4).Value = '''36.80'
$ws.Cells.Item(30, 4).Style = $origStyle
$origStyle = $ws.Cells.Item(30, 5).Style
$ws.Cells.Item(30, 5).Value = '''  +5.26%  '
$ws.Cells.Item(30, 5).Style = $origStyle

# Row 31
$origStyle = $ws.Cells.Item(31, 4).Style
$ws.Cells.Item(31, 4).Value = '''157.18'
$ws.Cells.Item(31, 4).Style = $origStyle
$origStyle = $ws.Cells.Item(31, 5).Style
$ws.Cells.Item(31, 5).Value = '''  +0.98%  '
$ws.Cells.Item(31, 5).Style = $origStyle

# Row 32
$origStyle = $ws.Cells.Item(32, 4).Style
$ws.Cells.Item(32, 4).Value = '''5.47'
$ws.Cells.Item(32, 4).Style = $origStyle
$origStyle = $ws.Cells.Item(32, 5).Style
$ws.Cells.Item(32, 5).Value = '''  +0.74%  '
$ws.Cells.Item(32, 5).Style = $origStyle

# Row 33
$origStyle = $ws.Cells.Item(33, 5).Style
$ws.Cells.Item(33, 5).Value = '''  -0.08%  '
$ws.Cells.Item(33, 5).Style = $origStyle

# Row 34
$origStyle = $ws.Cells.Item(34, 4).Style
$ws.Cells.Item(34, 4).Value = '''0.0765'
$ws.Cells.Item(34, 4).Style = $origStyle
$origStyle = $ws.Cells.Item(34, 5).Style
$ws.Cells.Item(34, 5).Value = '''  +0.91%  '
$ws.Cells.Item(34, 5).Style = $origStyle

# Row 35
$origStyle = $ws.Cells.Item(35, 5).Style
$ws.Cells.Item(35, 5).Value = '''  -0.55%  '
$ws.Cells.Item(35, 5).Style = $origStyle

# Row 36
$origStyle = $ws.Cells.Item(36, 4).Style
$ws.Cells.Item(36, 4).Value = '''17.32'
$ws.Cells.Item(36, 4).Style = $origStyle
$origStyle = $ws.Cells.Item(36, 5).Style
$ws.Cells.Item(36, 5).Value = '''  -0.85%  '
$ws.Cells.Item(36, 5).Style = $origStyle

# Row 37
$origStyle = $ws.Cells.Item(37, 5).Style
$ws.Cells.Item(37, 5).Value = '''  +4.52%  '
$ws.Cells.Item(37, 5).Style = $origStyle

# Row 38
$origStyle = $ws.Cells.Item(38, 5).Style
$ws.Cells.Item(38, 5).Value = '''  +1.87%  '
$ws.Cells.Item(38, 5).Style = $origStyle

# Row 39
$origStyle = $ws.Cells.Item(39, 5).Style
$ws.Cells.Item(39, 5).Value = '''  +0.82%  '
$ws.Cells.Item(39, 5).Style = $origStyle

# Row 40
$origStyle = $ws.Cells.Item(40, 5).Style
$ws.Cells.Item(40, 5).Value = '''  +2.46%  '
$ws.Cells.Item(40, 5).Style = $origStyle

# Row 41
$origStyle = $ws.Cells.Item(41, 5).Style
$ws.Cells.Item(41, 5).Value = '''  +1.07%  '
$ws.Cells.Item(41, 5).Style = $origStyle

# Row 42
$origStyle = $ws.Cells.Item(42, 4).Style
$ws.Cells.Item(42, 4).Value = '''2.49'
$ws.Cells.Item(42, 4).Style = $origStyle
$origStyle = $ws.Cells.Item(42, 5).Style
$ws.Cells.Item(42, 5).Value = '''  +2.65%  '
$ws.Cells.Item(42, 5).Style = $origStyle

# Row 43
$origStyle = $ws.Cells.Item(43, 4).Style
$ws.Cells.Item(43, 4).Value = '''1.986.33'
$ws.Cells.Item(43, 4).Style = $origStyle
$origStyle = $ws.Cells.Item(43, 5).Style
$ws.Cells.Item(43, 5).Value = '''  +1.18%  '
$ws.Cells.Item(43, 5).Style = $origStyle

# Row 44
$origStyle = $ws.Cells.Item(44, 2).Style
$ws.Cells.Item(44, 2).Value = '''EnergySwap'
$ws.Cells.Item(44, 2).Style = $origStyle
$origStyle = $ws.Cells.Item(44, 3).Style
$ws.Cells.Item(44, 3).Value = '''https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(44, 3).Style = $origStyle
$origStyle = $ws.Cells.Item(44, 4).Style
$ws.Cells.Item(44, 4).Value = '''18.90'
$ws.Cells.Item(44, 4).Style = $origStyle
$origStyle = $ws.Cells.Item(44, 5).Style
$ws.Cells.Item(44, 5).Value = '''  +2.33%  '
$ws.Cells.Item(44, 5).Style = $origStyle

# Row 45
$origStyle = $ws.Cells.Item(45, 2).Style
$ws.Cells.Item(45, 2).Value = '''VeChain'
$ws.Cells.Item(45, 2).Style = $origStyle
$origStyle = $ws.Cells.Item(45, 3).Style
$ws.Cells.Item(45, 3).Value = '''https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(45, 3).Style = $origStyle
$origStyle = $ws.Cells.Item(45, 4).Style
$ws.Cells.Item(45, 4).Value = '''0.0283'
$ws.Cells.Item(45, 4).Style = $origStyle
$origStyle = $ws.Cells.Item(45, 5).Style
$ws.Cells.Item(45, 5).Value = '''  +0.48%  '
$ws.Cells.Item(45, 5).Style = $origStyle

# Row 46
$origStyle = $ws.Cells.Item(46, 4).Style
$ws.Cells.Item(46, 4).Value = '''2.98'
$ws.Cells.Item(46, 4).Style = $origStyle
$origStyle = $ws.Cells.Item(46, 5).Style
$ws.Cells.Item(46, 5).Value = '''  +2.64%  '
$ws.Cells.Item(46, 5).Style = $origStyle

# Row 47
$origStyle = $ws.Cells.Item(47, 4).Style
$ws.Cells.Item(47, 4).Value = '''9.47'
$ws.Cells.Item(47, 4).Style = $origStyle
$origStyle = $ws.Cells.Item(47, 5).Style
$ws.Cells.Item(47, 5).Value = '''  +5.79%  '
$ws.Cells.Item(47, 5).Style = $origStyle

# Row 48
$origStyle = $ws.Cells.Item(48, 4).Style
$ws.Cells.Item(48, 4).Value = '''2.711.31'
$ws.Cells.Item(48, 4).Style = $origStyle
$origStyle = $ws.Cells.Item(48, 5).Style
$ws.Cells.Item(48, 5).Value = '''  +0.52%  '
$ws.Cells.Item(48, 5).Style = $origStyle

# Row 49
$origStyle = $ws.Cells.Item(49, 4).Style
$ws.Cells.Item(49, 4).Value = '''98.14'
$ws.Cells.Item(49, 4).Style = $origStyle
$origStyle = $ws.Cells.Item(49, 5).Style
$ws.Cells.Item(49, 5).Value = '''  +1.92%  '
$ws.Cells.Item(49, 5).Style = $origStyle

# Row 50
$origStyle = $ws.Cells.Item(50, 4).Style
$ws.Cells.Item(50, 4).Value = '''75.86'
$ws.Cells.Item(50, 4).Style = $origStyle
$origStyle = $ws.Cells.Item(50, 5).Style
$ws.Cells.Item(50, 5).Value = '''  +5.93%  '
$ws.Cells.Item(50, 5).Style = $origStyle

# Row 51
$origStyle = $ws.Cells.Item(51, 4).Style
$ws.Cells.Item(51, 4).Value = '''67.15'
$ws.Cells.Item(51, 4).Style = $origStyle
$origStyle = $ws.Cells.Item(51, 5).Style
$ws.Cells.Item(51, 5).Value = '''  +1.22%  '
$ws.Cells.Item(51, 5).Style = $origStyle
